# Fixing automation script for suite E
#
# The "Runmode" column (D) on the "Test Cases" sheet was wired to "N"
# (skip) for almost every scenario by mistake. Flip it back to "Y" so the
# automation actually runs suite E, and leave the sheet's selection on
# the range that was just fixed (matches what a user would see after
# multi-selecting + editing that column in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Column D (Runmode), rows 2-71: force every row back to "Y".
$ws.Range("D2:D71").Value = "Y"

# Reflect the edit in the sheet's selection/view state.
$ws.Range("D2:D71").Select()
